$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country name swaps (column A) driven by shared-string reordering ---
$ws.Range("A24").Value = "Colombia"
$ws.Range("A25").Value = "China"

$ws.Range("A72").Value = "Costa de Marfil"
$ws.Range("A73").Value = "Malasia"

$ws.Range("A200").Value = "Laos"
$ws.Range("A201").Value = "Santa Lucia"
$ws.Range("A202").Value = "Dominica"
$ws.Range("A203").Value = "Fiyi"

$ws.Range("A208").Value = "Islas Malvinas"
$ws.Range("A209").Value = "Groenlandia"

$ws.Range("A212").Value = "Montserrat"
$ws.Range("A213").Value = "Seychelles"

# --- Updated case numbers ---
$ws.Range("B4").Value = 2547365
$ws.Range("C4").Value = 42777
$ws.Range("D4").Value = 1060144
$ws.Range("E4").Value = 1359860
$ws.Range("G4").Value = 581
$ws.Range("H4").Value = 127361

$ws.Range("D5").Value = 697526
$ws.Range("E5").Value = 521487

$ws.Range("B17").Value = 194399
$ws.Range("C17").Value = 614
$ws.Range("E17").Value = 8273

$ws.Range("B22").Value = 102735
$ws.Range("C22").Value = 113
$ws.Range("D22").Value = 65658
$ws.Range("E22").Value = 28570

$ws.Range("B24").Value = 84442
$ws.Range("C24").Value = 3843
$ws.Range("D24").Value = 34937
$ws.Range("E24").Value = 46694
$ws.Range("G24").Value = 157
$ws.Range("H24").Value = 2811

$ws.Range("B25").Value = 83462
$ws.Range("C25").Value = 13
$ws.Range("D25").Value = 78439
$ws.Range("E25").Value = 389
$ws.Range("H25").Value = 4634

$ws.Range("D54").Value = 12548
$ws.Range("E54").Value = 7052
$ws.Range("G54").Value = 10
$ws.Range("H54").Value = 150

$ws.Range("B55").Value = 18197
$ws.Range("C55").Value = 87
$ws.Range("D55").Value = 16392
$ws.Range("E55").Value = 836
$ws.Range("G55").Value = 1
$ws.Range("H55").Value = 969

$ws.Range("B69").Value = 11038
$ws.Range("C69").Value = 168
$ws.Range("D69").Value = 7668
$ws.Range("E69").Value = 3021
$ws.Range("G69").Value = 4
$ws.Range("H69").Value = 349

$ws.Range("B70").Value = 9257
$ws.Range("C70").Value = 273
$ws.Range("D70").Value = 4014
$ws.Range("E70").Value = 4671
$ws.Range("G70").Value = 16
$ws.Range("H70").Value = 572

$ws.Range("B71").Value = 8832
$ws.Range("C71").Value = 44
$ws.Range("E71").Value = 445

$ws.Range("B72").Value = 8739
$ws.Range("C72").Value = 405
$ws.Range("D72").Value = 3587
$ws.Range("E72").Value = 5088
$ws.Range("G72").Value = 4
$ws.Range("H72").Value = 64

$ws.Range("B73").Value = 8606
$ws.Range("C73").Value = 6
$ws.Range("D73").Value = 8294
$ws.Range("E73").Value = 191
$ws.Range("H73").Value = 121

$ws.Range("B85").Value = 5260
$ws.Range("C85").Value = 86
$ws.Range("D85").Value = 4155
$ws.Range("E85").Value = 1076

$ws.Range("B86").Value = 5209
$ws.Range("C86").Value = 122
$ws.Range("D86").Value = 2327
$ws.Range("E86").Value = 2842

$ws.Range("B89").Value = 4513
$ws.Range("C89").Value = 105
$ws.Range("D89").Value = 2457
$ws.Range("E89").Value = 1841
$ws.Range("G89").Value = 4
$ws.Range("H89").Value = 215

$ws.Range("B104").Value = 2283
$ws.Range("C104").Value = 6
$ws.Range("E104").Value = 427

$ws.Range("B157").Value = 389
$ws.Range("C157").Value = 16
$ws.Range("D157").Value = 181
$ws.Range("E157").Value = 198

$ws.Range("B177").Value = 124
$ws.Range("C177").Value = 1
$ws.Range("E177").Value = 7

$ws.Range("D182").Value = 90
$ws.Range("E182").Value = 0

$ws.Range("D212").Value = 10
$ws.Range("H212").Value = 1

$ws.Range("D213").Value = 11
$ws.Range("H213").Value = 0

# --- Footer timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 27 de Junio de 2020 a las 00:50"
